$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update staffing numbers (column F) for Q2 2018 (rows 92-152) and
# --- the corresponding prior-year block (rows 273-303). Row 114 is left
# --- untouched, matching the source diff.

$values = @(69,94,127,130,126,66,68,101,115,78,100,109,140,113,120,58,136,85,89,134,96,69)
$row = 92
foreach ($v in $values) {
    $ws.Cells.Item($row, 6).Value = $v
    $row = $row + 1
}

$values = @(138,88,111,136,113,138,112,69,94,127,112,126,66,68,101,115,78,100,109,122,113,120,58,116,85,89,134,88,69,87,120,88,100,136,113,138,112,96)
$row = 115
foreach ($v in $values) {
    $ws.Cells.Item($row, 6).Value = $v
    $row = $row + 1
}

$values = @(60,82,109,112,109,59,63,82,103,25,89,96,122,100,102,51,113,77,81,118,84,63,73,111,48,97,123,95,98,116,96)
$row = 273
foreach ($v in $values) {
    $ws.Cells.Item($row, 6).Value = $v
    $row = $row + 1
}

# --- Reset the view: scroll back to the top of the sheet and select H1
# --- (previously the sheet had scrolled down to row 273 with H273 selected).
$ws.Range("H1").Select()
